$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("LaMelo Ball", "PG,SG", "Charlotte Hornets"),
    @("Brandon Miller", "SG,SF", "Charlotte Hornets"),
    @("Anthony Davis", "PF,C", "Los Angeles Lakers"),
    @("Robert Williams III", "C", "Portland Trail Blazers"),
    @("Bam Adebayo", "C", "Miami Heat"),
    @("T.J. McConnell", "PG", "Indiana Pacers"),
    @("Isaiah Hartenstein", "C", "Oklahoma City Thunder"),
    @("Julius Randle", "PF", "Minnesota Timberwolves"),
    @("Bilal Coulibaly", "SG,SF", "Washington Wizards"),
    @("Cam Thomas", "SG,SF", "Brooklyn Nets"),
    @("Jared McCain", "PG,SG", "Philadelphia 76ers"),
    @("Derrick White", "PG,SG", "Boston Celtics"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets")
)

$startRow = 3
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
